$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.255.62'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +3.74%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.331.93'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.03%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '547.22'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.51%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '131.40'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.07%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.581'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.50%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.329.35'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.02%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.101'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.12%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.52'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.98%  '
$ws.Range('E12').Value = '  +0.36%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.335'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.99%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.74'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.31%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.746.54'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.11%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '60.220.54'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.83%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000133'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.83%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.365.43'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.88%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.61'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.15%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.15'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.75%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '314.42'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.46%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.68'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.24%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.997'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.36%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '63.96'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.26%  '
$ws.Range('E25').Value = '  +1.59%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.84'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.51%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.36'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +7.83%  '
$ws.Range('B29').Value = 'SuiNetwork'
$ws.Range('C29').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.23'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +16.54%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '174.47'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.09%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.75'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.05%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0₃0729'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.92%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.95'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.42%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.38'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +11.70%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.380'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '17.87'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.48%  '
$ws.Range('E38').Value = '  +0.02%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.09'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +5.18%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '332.57'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +15.55%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '38.00'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.78%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.52'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.89%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '139.90'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.46%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.47'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.19%  '
$ws.Range('E45').Value = '  -0.91%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '19.39'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +8.29%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0496'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.72%  '
$ws.Range('E48').Value = '  +1.30%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0₆0230'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +22.01%  '
$ws.Range('E50').Value = '  +1.75%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '11.02'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.79%  '
